$d = $word.ActiveDocument

# The document contains a duplicated screenshot: the same "Extra - Tenacitas.png"
# picture appears once under "Introducción" (keep) and once again, redundantly,
# right under the "Presentación del proyecto" heading (remove). We locate the
# paragraph that (a) consists solely of a centered inline picture and (b) is
# immediately preceded by a paragraph whose text is "Presentación del proyecto",
# then delete that whole paragraph (including its paragraph mark), exactly as
# the target diff shows the <w:p> block being removed in full.

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $range = $para.Range

    if ($range.InlineShapes.Count -gt 0 -and $i -gt 1) {
        $prevText = $d.Paragraphs.Item($i - 1).Range.Text.Trim()
        if ($prevText -eq "Presentación del proyecto") {
            $range.Delete()
            break
        }
    }
}
